# Auto-generated PowerShell Excel COM-interop script
# Applies numeric cell value updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets
# (data refresh from a scheduled runner; see commit message)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3447.4
$ws.Range("I62").Value = 1925.75
$ws.Range("J62").Value = 4461.8335
$ws.Range("K62").Value = 1925.75
$ws.Range("L62").Value = 4461.8335
$ws.Range("M62").Value = -1301.75
$ws.Range("N62").Value = -5709.8335
$ws.Range("H65").Value = 3447.4
$ws.Range("I65").Value = 1925.75
$ws.Range("J65").Value = 4461.8335
$ws.Range("K65").Value = 9628.75
$ws.Range("L65").Value = 22309.1675
$ws.Range("M65").Value = -6508.75
$ws.Range("N65").Value = -28549.1675
$ws.Range("H76").Value = 3324.5833
$ws.Range("I76").Value = 3228
$ws.Range("K76").Value = 3228
$ws.Range("M76").Value = -2913
$ws.Range("H79").Value = 3324.5833
$ws.Range("I79").Value = 3228
$ws.Range("K79").Value = 3228
$ws.Range("M79").Value = -2136
$ws.Range("H106").Value = 918.1667
$ws.Range("I106").Value = 880.6429000000001
$ws.Range("J106").Value = 1049.5
$ws.Range("K106").Value = 880.6429000000001
$ws.Range("L106").Value = 1049.5
$ws.Range("M106").Value = -249.6429000000001
$ws.Range("N106").Value = -2311.5
$ws.Range("H132").Value = 10528056
$ws.Range("I132").Value = 11765945
$ws.Range("K132").Value = 35297835
$ws.Range("M132").Value = -35295305
$ws.Range("H137").Value = 7698912
$ws.Range("J137").Value = 2650.4443
$ws.Range("L137").Value = 7951.3329
$ws.Range("N137").Value = -13051.3329
$ws.Range("H141").Value = 581141.2
$ws.Range("I141").Value = 2323.75
$ws.Range("J141").Value = 1738776
$ws.Range("K141").Value = 6971.25
$ws.Range("L141").Value = 5216328
$ws.Range("M141").Value = -1791.25
$ws.Range("N141").Value = -5226688
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10002393
$ws.Range("I2").Value = 15626396
$ws.Range("J2").Value = 4166.6665
$ws.Range("K2").Value = 15626396
$ws.Range("L2").Value = 4166.6665
$ws.Range("M2").Value = -15626283
$ws.Range("N2").Value = -4392.6665
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H74").Value = 1403.4445
$ws.Range("I74").Value = 1252.6666
$ws.Range("K74").Value = 1252.6666
$ws.Range("M74").Value = -378.6666
$ws.Range("H77").Value = 1403.4445
$ws.Range("I77").Value = 1252.6666
$ws.Range("K77").Value = 6263.333000000001
$ws.Range("M77").Value = -1895.333000000001
$ws.Range("H80").Value = 30614
$ws.Range("J80").Value = 30614
$ws.Range("L80").Value = 30614
$ws.Range("N80").Value = -32610
$ws.Range("H83").Value = 30614
$ws.Range("J83").Value = 30614
$ws.Range("L83").Value = 91842
$ws.Range("N83").Value = -101826
$ws.Range("H97").Value = 487.17648
$ws.Range("I97").Value = 495.2857
$ws.Range("J97").Value = 449.33334
$ws.Range("K97").Value = 495.2857
$ws.Range("L97").Value = 449.33334
$ws.Range("M97").Value = 0.7142999999999802
$ws.Range("N97").Value = -1441.33334
$ws.Range("H102").Value = 2237.1428
$ws.Range("I102").Value = 1618.3334
$ws.Range("J102").Value = 5950
$ws.Range("K102").Value = 1618.3334
$ws.Range("L102").Value = 5950
$ws.Range("M102").Value = 3.666600000000017
$ws.Range("N102").Value = -9194
$ws.Range("H116").Value = 10002393
$ws.Range("I116").Value = 15626396
$ws.Range("J116").Value = 4166.6665
$ws.Range("K116").Value = 15626396
$ws.Range("L116").Value = 4166.6665
$ws.Range("M116").Value = -15624102
$ws.Range("N116").Value = -8754.666499999999
$ws.Range("H122").Value = 2507.9167
$ws.Range("I122").Value = 951
$ws.Range("J122").Value = 2819.3
$ws.Range("K122").Value = 2853
$ws.Range("L122").Value = 8457.900000000001
$ws.Range("M122").Value = -403
$ws.Range("N122").Value = -13357.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10002393
$ws.Range("I3").Value = 15626396
$ws.Range("J3").Value = 4166.6665
$ws.Range("K3").Value = 15626396
$ws.Range("L3").Value = 4166.6665
$ws.Range("M3").Value = -15626282
$ws.Range("N3").Value = -4394.6665
$ws.Range("H70").Value = 50384
$ws.Range("J70").Value = 50384
$ws.Range("L70").Value = 50384
$ws.Range("N70").Value = -50970
$ws.Range("H73").Value = 50384
$ws.Range("J73").Value = 50384
$ws.Range("L73").Value = 50384
$ws.Range("N73").Value = -52412
$ws.Range("H107").Value = 1637.1
$ws.Range("I107").Value = 1031.1818
$ws.Range("J107").Value = 2377.6667
$ws.Range("K107").Value = 1031.1818
$ws.Range("L107").Value = 2377.6667
$ws.Range("M107").Value = 888.8181999999999
$ws.Range("N107").Value = -6217.6667
$ws.Range("H134").Value = 2638.9285
$ws.Range("I134").Value = 1884.5
$ws.Range("J134").Value = 4525
$ws.Range("K134").Value = 5653.5
$ws.Range("L134").Value = 13575
$ws.Range("M134").Value = -3118.5
$ws.Range("N134").Value = -18645
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2781024.8
$ws.Range("I31").Value = 3450127.2
$ws.Range("J31").Value = 9028.571
$ws.Range("K31").Value = 3450127.2
$ws.Range("L31").Value = 9028.571
$ws.Range("M31").Value = -3449832.2
$ws.Range("N31").Value = -9618.571
$ws.Range("H33").Value = 20426.637
$ws.Range("I33").Value = 1823.5
$ws.Range("K33").Value = 1823.5
$ws.Range("M33").Value = -1444.5
$ws.Range("H34").Value = 2781024.8
$ws.Range("I34").Value = 3450127.2
$ws.Range("J34").Value = 9028.571
$ws.Range("K34").Value = 3450127.2
$ws.Range("L34").Value = 9028.571
$ws.Range("M34").Value = -3449925.2
$ws.Range("N34").Value = -9432.571
$ws.Range("H36").Value = 38684.332
$ws.Range("I36").Value = 21250
$ws.Range("J36").Value = 52631.8
$ws.Range("K36").Value = 21250
$ws.Range("L36").Value = 52631.8
$ws.Range("M36").Value = -20862
$ws.Range("N36").Value = -53407.8
$ws.Range("H40").Value = 38684.332
$ws.Range("I40").Value = 21250
$ws.Range("J40").Value = 52631.8
$ws.Range("K40").Value = 21250
$ws.Range("L40").Value = 52631.8
$ws.Range("M40").Value = -21090
$ws.Range("N40").Value = -52951.8
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H58").Value = 15627580
$ws.Range("I58").Value = 1473.45
$ws.Range("J58").Value = 41671092
$ws.Range("K58").Value = 1473.45
$ws.Range("L58").Value = 41671092
$ws.Range("M58").Value = -1270.45
$ws.Range("N58").Value = -41671498
$ws.Range("H94").Value = 23811194
$ws.Range("I94").Value = 1401
$ws.Range("J94").Value = 33335112
$ws.Range("K94").Value = 1401
$ws.Range("L94").Value = 33335112
$ws.Range("M94").Value = -950
$ws.Range("N94").Value = -33336014
$ws.Range("H132").Value = 4027.76
$ws.Range("I132").Value = 1917.091
$ws.Range("J132").Value = 5686.143
$ws.Range("K132").Value = 5751.272999999999
$ws.Range("L132").Value = 17058.429
$ws.Range("M132").Value = -3221.272999999999
$ws.Range("N132").Value = -22118.429
$ws.Range("H134").Value = 2326.8572
$ws.Range("I134").Value = 1770.9333
$ws.Range("J134").Value = 3716.6667
$ws.Range("K134").Value = 5312.7999
$ws.Range("L134").Value = 11150.0001
$ws.Range("M134").Value = -2777.7999
$ws.Range("N134").Value = -16220.0001
$ws.Range("H136").Value = 15627580
$ws.Range("I136").Value = 1473.45
$ws.Range("J136").Value = 41671092
$ws.Range("K136").Value = 4420.35
$ws.Range("L136").Value = 125013276
$ws.Range("M136").Value = -1870.35
$ws.Range("N136").Value = -125018376
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 186.66667
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 186.66667
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 560.00001
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1030.00001
$ws.Range("H107").Value = 1409.091
$ws.Range("I107").Value = 940
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 2820
$ws.Range("L107").Value = 5400
$ws.Range("M107").Value = -900
$ws.Range("N107").Value = -9240
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4284.615
$ws.Range("I80").Value = 3916.6667
$ws.Range("J80").Value = 4600
$ws.Range("K80").Value = 3916.6667
$ws.Range("L80").Value = 4600
$ws.Range("M80").Value = -2918.6667
$ws.Range("N80").Value = -6596
$ws.Range("H83").Value = 4284.615
$ws.Range("I83").Value = 3916.6667
$ws.Range("J83").Value = 4600
$ws.Range("K83").Value = 19583.3335
$ws.Range("L83").Value = 23000
$ws.Range("M83").Value = -14591.3335
$ws.Range("N83").Value = -32984
$ws.Range("H102").Value = 2040.409
$ws.Range("I102").Value = 1592.0714
$ws.Range("J102").Value = 2825
$ws.Range("K102").Value = 1592.0714
$ws.Range("L102").Value = 2825
$ws.Range("M102").Value = 29.92859999999996
$ws.Range("N102").Value = -6069
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2198.8823
$ws.Range("I7").Value = 1658.1
$ws.Range("K7").Value = 1658.1
$ws.Range("M7").Value = -1546.1
$ws.Range("H126").Value = 2198.8823
$ws.Range("I126").Value = 1658.1
$ws.Range("K126").Value = 4974.299999999999
$ws.Range("M126").Value = -2504.299999999999
